# "Rcc script is completed" - add the RCC011 test case row to the
# "Test Cases" sheet (mirrors the existing rows, which already follow the
# TCID / Jira id / Description / Runmode / Results layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Copy the formatting of the last existing data row (row 9) down onto the
# new row 10, the same way a user extending this list in Excel would do,
# then overwrite it with the new test case's values.
$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = "RCC011"
$ws.Range("B10").Value = "OBT2"
$ws.Range("C10").Value = "Verify the invitations scenario"
$ws.Range("D10").Value = "Y"

# Leave the selection where it ended up after finishing the entry.
$ws.Range("C14").Select() | Out-Null
